function SetTextValue($Range, $Text) {
    # Write a value while keeping the cell a text/string cell (matches
    # the source workbook, where numeric-looking values like prices and
    # hour counters are stored as text, e.g. t="inlineStr").
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

SetTextValue $ws.Range("D2") "244.09"
SetTextValue $ws.Range("G2") "5"
SetTextValue $ws.Range("D3") "23.09"
SetTextValue $ws.Range("G3") "5"
SetTextValue $ws.Range("G4") "5"
SetTextValue $ws.Range("D5") "0.05974"
SetTextValue $ws.Range("G5") "5"
SetTextValue $ws.Range("D6") "3.458"
SetTextValue $ws.Range("G6") "5"
SetTextValue $ws.Range("D7") "6.529"
SetTextValue $ws.Range("G7") "5"
SetTextValue $ws.Range("D8") "0.8159"
SetTextValue $ws.Range("G8") "5"
SetTextValue $ws.Range("D9") "0.9150"
SetTextValue $ws.Range("G9") "5"
SetTextValue $ws.Range("G10") "5"
SetTextValue $ws.Range("D11") "0.07415"
SetTextValue $ws.Range("G11") "5"
SetTextValue $ws.Range("D12") "0.03253"
SetTextValue $ws.Range("G12") "5"
SetTextValue $ws.Range("D13") "0.03066"
SetTextValue $ws.Range("G13") "5"
SetTextValue $ws.Range("G14") "5"
SetTextValue $ws.Range("D15") "3.861"
SetTextValue $ws.Range("G15") "5"
SetTextValue $ws.Range("D16") "0.001577"
SetTextValue $ws.Range("G16") "5"
SetTextValue $ws.Range("D17") "0.04662"
SetTextValue $ws.Range("G17") "5"
SetTextValue $ws.Range("D18") "0.0005940"
SetTextValue $ws.Range("G18") "5"
SetTextValue $ws.Range("D19") "0.006095"
SetTextValue $ws.Range("G19") "5"
SetTextValue $ws.Range("D20") "0.005004"
SetTextValue $ws.Range("G20") "5"
SetTextValue $ws.Range("D21") "0.0009860"
SetTextValue $ws.Range("G21") "5"
SetTextValue $ws.Range("G22") "5"
SetTextValue $ws.Range("D23") "3.613"
SetTextValue $ws.Range("G23") "5"
SetTextValue $ws.Range("D24") "2.138"
SetTextValue $ws.Range("G24") "5"
SetTextValue $ws.Range("G25") "5"
SetTextValue $ws.Range("D26") "0.1296"
SetTextValue $ws.Range("G26") "5"
SetTextValue $ws.Range("D27") "0.0002394"
SetTextValue $ws.Range("G27") "5"
SetTextValue $ws.Range("G28") "5"
SetTextValue $ws.Range("G29") "5"
SetTextValue $ws.Range("G30") "5"
SetTextValue $ws.Range("G31") "5"
SetTextValue $ws.Range("G32") "5"
SetTextValue $ws.Range("G33") "5"
SetTextValue $ws.Range("G34") "5"
SetTextValue $ws.Range("G35") "5"
SetTextValue $ws.Range("G36") "5"
SetTextValue $ws.Range("G37") "5"
SetTextValue $ws.Range("G38") "5"
SetTextValue $ws.Range("G39") "5"
SetTextValue $ws.Range("D40") "0.03933"
SetTextValue $ws.Range("G40") "5"
SetTextValue $ws.Range("D41") "0.006207"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
SetTextValue $ws.Range("G41") "5"
SetTextValue $ws.Range("D42") "0.1074"
SetTextValue $ws.Range("G42") "5"
SetTextValue $ws.Range("D43") "0.003200"
SetTextValue $ws.Range("G43") "5"
SetTextValue $ws.Range("D44") "0.007242"
$ws.Range("E44").Value = "43LocalTradersLCT"
SetTextValue $ws.Range("G44") "5"
SetTextValue $ws.Range("D45") "0.00005233"
SetTextValue $ws.Range("G45") "5"
SetTextValue $ws.Range("G46") "5"
SetTextValue $ws.Range("G47") "5"
SetTextValue $ws.Range("D48") "0.7800"
SetTextValue $ws.Range("G48") "5"
SetTextValue $ws.Range("G49") "5"
SetTextValue $ws.Range("G50") "5"
SetTextValue $ws.Range("G51") "5"
